$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.808.51"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "3.758.55"
$ws.Range("E3").Value = "  -1.58%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'635.69"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").Value = "'165.30"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").Value = "3.757.74"
$ws.Range("E7").Value = "  -1.53%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  -2.70%  "

$ws.Range("D11").Value = "'0.456"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("E12").Value = "  +4.40%  "

$ws.Range("D13").Value = "'0.0000238"
$ws.Range("E13").Value = "  -5.09%  "

$ws.Range("D14").Value = "'34.81"
$ws.Range("E14").Value = "  -3.41%  "

$ws.Range("D15").Value = "4.390.49"
$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("D16").Value = "3.760.83"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").Value = "68.816.96"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").Value = "'17.60"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").Value = "'0.114"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").Value = "'6.98"
$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("D21").Value = "'469.29"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "'9.49"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").Value = "'81.50"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("E25").Value = "  -6.75%  "

$ws.Range("D26").Value = "'12.12"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").Value = "'2.10"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "3.906.40"
$ws.Range("E30").Value = "  -1.50%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.68"
$ws.Range("E31").Value = "  -1.09%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.25"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").Value = "'7.09"
$ws.Range("E33").Value = "  -2.90%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'28.38"
$ws.Range("E34").Value = "  -2.93%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.174"
$ws.Range("E35").Value = "  +16.15%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "3.711.62"
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("D38").Value = "'8.85"
$ws.Range("E38").Value = "  -3.01%  "

$ws.Range("E39").Value = "  -1.43%  "

$ws.Range("E40").Value = "  -5.77%  "

$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  -2.88%  "

$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "'0.954"
$ws.Range("E43").Value = "  -2.92%  "

$ws.Range("D45").Value = "'44.43"
$ws.Range("E45").Value = "  +3.71%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.97"
$ws.Range("E46").Value = "  +3.41%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'155.40"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("D48").Value = "'47.39"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("E49").Value = "  -3.54%  "

$ws.Range("E50").Value = "  -2.77%  "

$ws.Range("D51").Value = "'8.34"
$ws.Range("E51").Value = "  -1.38%  "
